$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New weekly section (week of 2018.10.10), mirroring the layout of the
# previous weeks already present in rows 1-18.

# Date header row
$ws.Range("A19").Value = "日期：2018.10.10 第六周周三"

# Column headers row
$ws.Range("A20").Value = "组员"
$ws.Range("B20").Value = "计划内容"
$ws.Range("C20").Value = "完成情况"
$ws.Range("D20").Value = "备注"

# Team member rows
$ws.Range("A21").Value = "余舒章"
$ws.Range("B21").Value = "写6条用例规约"
$ws.Range("C21").Value = "已完成"

$ws.Range("A22").Value = "王嘉宇"
$ws.Range("B22").Value = "写6条用例规约"
$ws.Range("C22").Value = "已完成"

$ws.Range("A23").Value = "许俊杰"
$ws.Range("B23").Value = "写6条用例规约"
$ws.Range("C23").Value = "已完成"

$ws.Range("A24").Value = "庞森杰"
$ws.Range("B24").Value = "写6条用例规约"
$ws.Range("C24").Value = "已完成"

# Summary row (row 25 left blank, matching the blank spacer row pattern)
$ws.Range("A26").Value = "总结："

# Update selection to match the author's final cursor position
[void]$ws.Range("C13").Select()
